$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidentiality / as-of date note in cell A18
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-04-05 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15

$ws.Range("D2").Value = 0.05640035411919588
$ws.Range("E2").Value = 0.01165902140672781

$ws.Range("D3").Value = 0.02346309169773744
$ws.Range("E3").Value = 0.01507537688442206

$ws.Range("D4").Value = 0.03184647288926499
$ws.Range("E4").Value = 0.002094440213252113

$ws.Range("D5").Value = 0.03211799122247929
$ws.Range("E5").Value = -0.01195695496213633

$ws.Range("D6").Value = 0.03440130544397627
$ws.Range("E6").Value = 0.02397179788484149

$ws.Range("D7").Value = 0.01893082413562172
$ws.Range("E7").Value = 0.005119581464872791

$ws.Range("D8").Value = 0.004887329997857498
$ws.Range("E8").Value = -0.01530190239867657

$ws.Range("D9").Value = 0.006721257796888118
$ws.Range("E9").Value = 0.01242983159583

$ws.Range("D10").Value = 0.06944940393977818
$ws.Range("E10").Value = -0.001164144353899887

$ws.Range("D11").Value = 0.0695706776369955
$ws.Range("E11").Value = -0.0005810575246950611

$ws.Range("D12").Value = 0.1482341875942398
$ws.Range("E12").Value = -0.004363319031343127

$ws.Range("D13").Value = 0.3895149456356965
$ws.Range("E13").Value = -0.002822865208186309

$ws.Range("D14").Value = 0.1144621578902689
$ws.Range("E14").Value = -0.004820766378244801

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -0.0007951175209500727

# Restore sheet protection (as it was before the edit)
$ws.Protect($null, $true, $true, $true)
